$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.930.31'
$ws.Range('E2').Value = '  +2.65%  '
$ws.Range('D3').Value = '2.229.49'
$ws.Range('E3').Value = '  +1.06%  '
$ws.Range('E4').Value = '  -0.03%  '
$c = $ws.Range('D5')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '293.03'
$c.Style = $origStyle
$ws.Range('E5').Value = '  -0.92%  '
$c = $ws.Range('D6')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '85.80'
$c.Style = $origStyle
$ws.Range('E6').Value = '  +5.03%  '
$c = $ws.Range('D7')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.515'
$c.Style = $origStyle
$ws.Range('E7').Value = '  +1.59%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +1.56%  '
$c = $ws.Range('D10')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '30.69'
$c.Style = $origStyle
$ws.Range('E10').Value = '  +7.94%  '
$ws.Range('E11').Value = '  +2.13%  '
$c = $ws.Range('D12')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '47.04'
$c.Style = $origStyle
$ws.Range('E12').Value = '  +0.56%  '
$c = $ws.Range('D13')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.109'
$c.Style = $origStyle
$ws.Range('E13').Value = '  +1.38%  '
$c = $ws.Range('D14')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '6.39'
$c.Style = $origStyle
$ws.Range('E14').Value = '  +3.44%  '
$ws.Range('D15').Value = '2.572.65'
$ws.Range('E15').Value = '  +0.99%  '
$c = $ws.Range('D16')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '14.08'
$c.Style = $origStyle
$ws.Range('E16').Value = '  +0.62%  '
$ws.Range('D17').Value = '2.227.55'
$ws.Range('E17').Value = '  +0.36%  '
$c = $ws.Range('D18')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.727'
$c.Style = $origStyle
$ws.Range('E18').Value = '  +2.25%  '
$ws.Range('D19').Value = '39.867.12'
$ws.Range('E19').Value = '  +2.72%  '
$ws.Range('E20').Value = '  +2.65%  '
$ws.Range('E21').Value = '  +1.72%  '
$c = $ws.Range('D22')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '10.76'
$c.Style = $origStyle
$ws.Range('E22').Value = '  +7.46%  '
$c = $ws.Range('D23')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '65.04'
$c.Style = $origStyle
$ws.Range('E23').Value = '  +0.80%  '
$c = $ws.Range('D24')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '235.21'
$c.Style = $origStyle
$ws.Range('E24').Value = '  +4.55%  '
$c = $ws.Range('D25')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = $origStyle
$ws.Range('E25').Value = '  -0.15%  '
$c = $ws.Range('D26')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.46'
$c.Style = $origStyle
$ws.Range('E26').Value = '  +3.41%  '
$ws.Range('E27').Value = '  +5.30%  '
$c = $ws.Range('D28')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '22.77'
$c.Style = $origStyle
$ws.Range('E28').Value = '  +1.75%  '
$c = $ws.Range('D29')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.22'
$c.Style = $origStyle
$ws.Range('E29').Value = '  +2.96%  '
$c = $ws.Range('D30')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '9.24'
$c.Style = $origStyle
$ws.Range('E30').Value = '  +2.86%  '
$c = $ws.Range('D31')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '32.98'
$c.Style = $origStyle
$ws.Range('E31').Value = '  +4.95%  '
$c = $ws.Range('D32')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '152.19'
$c.Style = $origStyle
$ws.Range('E32').Value = '  +3.03%  '
$c = $ws.Range('D33')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = $origStyle
$ws.Range('E33').Value = '  -0.08%  '
$c = $ws.Range('D34')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.90'
$c.Style = $origStyle
$ws.Range('E34').Value = '  +1.78%  '
$c = $ws.Range('D35')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0718'
$c.Style = $origStyle
$ws.Range('E35').Value = '  +3.89%  '
$ws.Range('E36').Value = '  +3.38%  '
$c = $ws.Range('D37')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '16.16'
$c.Style = $origStyle
$ws.Range('E37').Value = '  +11.12%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D38')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.79'
$c.Style = $origStyle
$ws.Range('E38').Value = '  +5.78%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D39')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.111'
$c.Style = $origStyle
$ws.Range('E39').Value = '  +1.72%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D40')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0996'
$c.Style = $origStyle
$ws.Range('E40').Value = '  +5.22%  '
$ws.Range('E41').Value = '  +5.76%  '
$c = $ws.Range('D42')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.82'
$c.Style = $origStyle
$ws.Range('E42').Value = '  +4.60%  '
$ws.Range('D43').Value = '2.039.94'
$ws.Range('E43').Value = '  +7.64%  '
$ws.Range('E44').Value = '  +7.86%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D45')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0270'
$c.Style = $origStyle
$ws.Range('E45').Value = '  +5.08%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D46')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '10.00'
$c.Style = $origStyle
$ws.Range('E46').Value = '  +11.78%  '
$c = $ws.Range('D47')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '16.60'
$c.Style = $origStyle
$ws.Range('E47').Value = '  +4.11%  '
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('D49').Value = '2.444.75'
$ws.Range('E49').Value = '  +1.60%  '
$c = $ws.Range('D50')
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '70.59'
$c.Style = $origStyle
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('E51').Value = '  +2.96%  '
